$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.899.46"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.668.09"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'215.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("E6").Value = "  +5.55%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.0620"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("D9").Value = "'0.250"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "'20.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("E11").Value = "  +3.63%  "
$ws.Range("D12").Value = "1.902.58"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "1.654.80"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("D17").Value = "26.926.57"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'234.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").Value = "0.0₃0732"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "'2.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("D24").Value = "'9.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Value = "'146.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").Value = "'7.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("D28").Value = "'15.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("D33").Value = "1.449.39"
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("D34").Value = "'3.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("D35").Value = "'1.66"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.24%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").Value = "'0.581"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "'0.904"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").Value = "'66.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").Value = "'0.973"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.96%  "
$ws.Range("D45").Value = "1.811.06"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").Value = "'0.783"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").Value = "'90.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("E49").Value = "  +4.89%  "
$ws.Range("D50").Value = "'0.0505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").Value = "'7.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.57%  "
